# Atualização de bases das ligas, do dia: 15-06-2024 às 21:10
# Swap the match-detail values (columns B through AD) between pairs of
# adjacent rows in the active sheet. The row-index column (A) is left
# untouched; only the data describing each match (id, teams, scores,
# odds, etc.) moves between the two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 2   # column B
$lastCol  = 30  # column AD

$rowPairs = @(
    @(61, 62),
    @(156, 157),
    @(187, 188),
    @(228, 229),
    @(252, 253)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell1 = $ws.Cells.Item($r1, $c)
        $cell2 = $ws.Cells.Item($r2, $c)

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}
